$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting rows 93:165 down to 94:166
$ws.Rows.Item(93).Insert()

# Fill in the new row 93 with the new record's data
$ws.Range("A93").Value = 11
$ws.Range("B93").Value = "Vega Monumental Concepción"
$ws.Range("C93").Value = "Bíobío"
$ws.Range("D93").Value = 44957
$ws.Range("D93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E93").Value = 8
$ws.Range("F93").Value = 100112021
$ws.Range("G93").Value = "Ají"
$ws.Range("H93").Value = "Americana (o)"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 35
$ws.Range("K93").Value = 18000
$ws.Range("L93").Value = 20000
$ws.Range("M93").Value = 19143
$ws.Range("N93").Value = "$/caja 10 kilos"
$ws.Range("O93").Value = "Región Metropolitana"
$ws.Range("P93").Value = 1914
$ws.Range("Q93").Value = 10
$ws.Range("R93").Value = "Hortaliza"
